$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.293826460838318
$ws.Range("B1").Value = 3.124848365783691
$ws.Range("C1").Value = 2.489050149917603
$ws.Range("D1").Value = 2.341967344284058
$ws.Range("E1").Value = 1.989774942398071
